$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 210.55173
$ws.Range("I33").Value = 154
$ws.Range("K33").Value = 154
$ws.Range("M33").Value = 75
$ws.Range("H121").Value = 1299
$ws.Range("I121").Value = 595
$ws.Range("J121").Value = 1475
$ws.Range("K121").Value = 1785
$ws.Range("L121").Value = 4425
$ws.Range("M121").Value = -38
$ws.Range("N121").Value = -7919
$ws.Range("H127").Value = 1386.25
$ws.Range("I127").Value = 829.375
$ws.Range("J127").Value = 2500
$ws.Range("K127").Value = 2488.125
$ws.Range("L127").Value = 7500
$ws.Range("M127").Value = 2471.875
$ws.Range("N127").Value = -17420
$ws.Range("H132").Value = 360505.66
$ws.Range("I132").Value = 373813.25
$ws.Range("K132").Value = 1121439.75
$ws.Range("M132").Value = -1118909.75
$ws.Range("H137").Value = 67664456
$ws.Range("I137").Value = 1660085.6
$ws.Range("J137").Value = 166671000
$ws.Range("K137").Value = 4980256.800000001
$ws.Range("L137").Value = 500013000
$ws.Range("M137").Value = -4977706.800000001
$ws.Range("N137").Value = -500018100
$ws.Range("H141").Value = 1422.5
$ws.Range("I141").Value = 1422.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4267.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 912.5
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4912.6445
$ws.Range("I32").Value = 5199.468
$ws.Range("J32").Value = 3642.4285
$ws.Range("K32").Value = 5199.468
$ws.Range("L32").Value = 3642.4285
$ws.Range("M32").Value = -4912.468
$ws.Range("N32").Value = -4216.4285
$ws.Range("H61").Value = 2874.6296
$ws.Range("I61").Value = 2288.1052
$ws.Range("J61").Value = 4267.625
$ws.Range("K61").Value = 2288.1052
$ws.Range("L61").Value = 4267.625
$ws.Range("M61").Value = -2076.1052
$ws.Range("N61").Value = -4691.625
$ws.Range("H74").Value = 5266.357
$ws.Range("I74").Value = 1528.8
$ws.Range("J74").Value = 6078.8696
$ws.Range("K74").Value = 1528.8
$ws.Range("L74").Value = 6078.8696
$ws.Range("M74").Value = -654.8
$ws.Range("N74").Value = -7826.8696
$ws.Range("H77").Value = 5266.357
$ws.Range("I77").Value = 1528.8
$ws.Range("J77").Value = 6078.8696
$ws.Range("K77").Value = 7644
$ws.Range("L77").Value = 30394.348
$ws.Range("M77").Value = -3276
$ws.Range("N77").Value = -39130.348
$ws.Range("H88").Value = 2750.8125
$ws.Range("I88").Value = 2884
$ws.Range("J88").Value = 2351.25
$ws.Range("K88").Value = 2884
$ws.Range("L88").Value = 2351.25
$ws.Range("M88").Value = -2478
$ws.Range("N88").Value = -3163.25
$ws.Range("H91").Value = 2750.8125
$ws.Range("I91").Value = 2884
$ws.Range("J91").Value = 2351.25
$ws.Range("K91").Value = 2884
$ws.Range("L91").Value = 2351.25
$ws.Range("M91").Value = -1480
$ws.Range("N91").Value = -5159.25
$ws.Range("H122").Value = 2001.037
$ws.Range("I122").Value = 1718.6522
$ws.Range("K122").Value = 5155.9566
$ws.Range("M122").Value = -2705.9566
$ws.Range("H132").Value = 1328579.2
$ws.Range("I132").Value = 1817345.4
$ws.Range("K132").Value = 5452036.199999999
$ws.Range("M132").Value = -5449506.199999999
$ws.Range("H136").Value = 2874.6296
$ws.Range("I136").Value = 2288.1052
$ws.Range("J136").Value = 4267.625
$ws.Range("K136").Value = 6864.3156
$ws.Range("L136").Value = 12802.875
$ws.Range("M136").Value = -4314.3156
$ws.Range("N136").Value = -17902.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1604.5814
$ws.Range("I86").Value = 1490.9706
$ws.Range("J86").Value = 2033.7778
$ws.Range("K86").Value = 1490.9706
$ws.Range("L86").Value = 2033.7778
$ws.Range("M86").Value = -367.9706000000001
$ws.Range("N86").Value = -4279.7778
$ws.Range("H89").Value = 1604.5814
$ws.Range("I89").Value = 1490.9706
$ws.Range("J89").Value = 2033.7778
$ws.Range("K89").Value = 7454.853000000001
$ws.Range("L89").Value = 10168.889
$ws.Range("M89").Value = -1838.853000000001
$ws.Range("N89").Value = -21400.889
$ws.Range("H134").Value = 54185.43
$ws.Range("I134").Value = 56744.7
$ws.Range("K134").Value = 170234.1
$ws.Range("M134").Value = -167699.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1329.8
$ws.Range("I31").Value = 990.4815
$ws.Range("J31").Value = 2475
$ws.Range("K31").Value = 990.4815
$ws.Range("L31").Value = 2475
$ws.Range("M31").Value = -695.4815
$ws.Range("N31").Value = -3065
$ws.Range("H34").Value = 1329.8
$ws.Range("I34").Value = 990.4815
$ws.Range("J34").Value = 2475
$ws.Range("K34").Value = 990.4815
$ws.Range("L34").Value = 2475
$ws.Range("M34").Value = -788.4815
$ws.Range("N34").Value = -2879
$ws.Range("H58").Value = 3173.9023
$ws.Range("I58").Value = 1642.75
$ws.Range("J58").Value = 4153.84
$ws.Range("K58").Value = 1642.75
$ws.Range("L58").Value = 4153.84
$ws.Range("M58").Value = -1439.75
$ws.Range("N58").Value = -4559.84
$ws.Range("H97").Value = 13900
$ws.Range("J97").Value = 13900
$ws.Range("L97").Value = 13900
$ws.Range("N97").Value = -15882
$ws.Range("H99").Value = 1056
$ws.Range("I99").Value = 1056
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1056
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 442
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 1056
$ws.Range("I126").Value = 1056
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3168
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -698
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2105.3096
$ws.Range("I132").Value = 1911.9259
$ws.Range("J132").Value = 2453.4
$ws.Range("K132").Value = 5735.7777
$ws.Range("L132").Value = 7360.200000000001
$ws.Range("M132").Value = -3205.7777
$ws.Range("N132").Value = -12420.2
$ws.Range("H134").Value = 2457.4644
$ws.Range("I134").Value = 2546.7727
$ws.Range("J134").Value = 2130
$ws.Range("K134").Value = 7640.3181
$ws.Range("L134").Value = 6390
$ws.Range("M134").Value = -5105.3181
$ws.Range("N134").Value = -11460
$ws.Range("H136").Value = 3173.9023
$ws.Range("I136").Value = 1642.75
$ws.Range("J136").Value = 4153.84
$ws.Range("K136").Value = 4928.25
$ws.Range("L136").Value = 12461.52
$ws.Range("M136").Value = -2378.25
$ws.Range("N136").Value = -17561.52

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2810.682
$ws.Range("I113").Value = 520.625
$ws.Range("J113").Value = 5558.75
$ws.Range("K113").Value = 1561.875
$ws.Range("L113").Value = 16676.25
$ws.Range("M113").Value = 608.125
$ws.Range("N113").Value = -21016.25
$ws.Range("H120").Value = 5936.8
$ws.Range("I120").Value = 5684
$ws.Range("J120").Value = 6000
$ws.Range("K120").Value = 17052
$ws.Range("L120").Value = 18000
$ws.Range("M120").Value = -12214
$ws.Range("N120").Value = -27676
$ws.Range("H131").Value = 1726593.1
$ws.Range("J131").Value = 2327629.8
$ws.Range("L131").Value = 6982889.399999999
$ws.Range("N131").Value = -6992969.399999999
$ws.Range("H133").Value = 6926.8
$ws.Range("J133").Value = 7999.067
$ws.Range("L133").Value = 23997.201
$ws.Range("N133").Value = -34117.201

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 464153.84
$ws.Range("J21").Value = 2001333.4
$ws.Range("L21").Value = 2001333.4
$ws.Range("N21").Value = -2001679.4
$ws.Range("H30").Value = 464153.84
$ws.Range("J30").Value = 2001333.4
$ws.Range("L30").Value = 2001333.4
$ws.Range("N30").Value = -2001543.4
$ws.Range("H116").Value = 54785
$ws.Range("J116").Value = 54785
$ws.Range("L116").Value = 54785
$ws.Range("N116").Value = -63963
$ws.Range("H132").Value = 3140
$ws.Range("I132").Value = 2444.1333
$ws.Range("J132").Value = 5749.5
$ws.Range("K132").Value = 7332.3999
$ws.Range("L132").Value = 17248.5
$ws.Range("M132").Value = -4802.3999
$ws.Range("N132").Value = -22308.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7178
$ws.Range("I132").Value = 11542.667
$ws.Range("K132").Value = 34628.001
$ws.Range("M132").Value = -32098.001
$ws.Range("H136").Value = 1683.125
$ws.Range("I136").Value = 780
$ws.Range("J136").Value = 3188.3333
$ws.Range("K136").Value = 2340
$ws.Range("L136").Value = 9564.999899999999
$ws.Range("M136").Value = 210
$ws.Range("N136").Value = -14664.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4797.049
$ws.Range("I132").Value = 5652.129
$ws.Range("K132").Value = 16956.387
$ws.Range("M132").Value = -14426.387
$ws.Range("H136").Value = 10156.454
$ws.Range("I136").Value = 15132.643
$ws.Range("J136").Value = 1448.125
$ws.Range("K136").Value = 45397.929
$ws.Range("L136").Value = 4344.375
$ws.Range("M136").Value = -42847.929
$ws.Range("N136").Value = -9444.375
